$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.166.77'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '1.840.93'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.76'
$ws.Range('E5').Value = '  -3.04%  '
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4640'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3861'
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07854'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9623'
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.06'
$ws.Range('E11').Value = '  -0.82%  '
$ws.Range('D12').Value = '1.834.83'
$ws.Range('E12').Value = '  -3.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.691'
$ws.Range('E13').Value = '  -2.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.874'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06902'
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.52'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009955'
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.70'
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').Value = '28.187.60'
$ws.Range('E21').Value = '  -1.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.304'
$ws.Range('E22').Value = '  -1.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.02'
$ws.Range('E23').Value = '  -1.78%  '
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('D25').Value = '2.078.89'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.49'
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.16'
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.742'
$ws.Range('E28').Value = '  -5.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.968'
$ws.Range('E29').Value = '  -1.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.76'
$ws.Range('E30').Value = '  +0.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09262'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9283'
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.283'
$ws.Range('E33').Value = '  -1.54%  '
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.332'
$ws.Range('E35').Value = '  -3.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05805'
$ws.Range('E36').Value = '  -4.92%  '
$ws.Range('E37').Value = '  -4.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.145'
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.748'
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.885'
$ws.Range('E41').Value = '  -2.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1762'
$ws.Range('E42').Value = '  -1.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07223'
$ws.Range('E43').Value = '  +1.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.58'
$ws.Range('E44').Value = '  -1.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5269'
$ws.Range('E45').Value = '  -1.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.132'
$ws.Range('E46').Value = '  -9.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.117'
$ws.Range('E47').Value = '  -12.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.831'
$ws.Range('E48').Value = '  -3.64%  '
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('B51').Value = 'MXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.323'
$ws.Range('E51').Value = '  -0.86%  '
